$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from existing header cell (e.g. AC1) to new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 69   # AD
    $ws.Cells.Item($r, 31).Value = 93   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
